$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New labels / text cells, entered in the order that reproduces the
# --- shared-string table order seen in the target workbook.
$ws.Range("A18").Value = "Max size"
$ws.Range("D20").Value = "works"
$ws.Range("B19").Value = "height"
$ws.Range("F19").Value = "total pixels "
$ws.Range("A22").Value = "try"
$ws.Range("C19").Value = "wide"
$ws.Range("H19").Value = "notes"
$ws.Range("H20").Value = "15 inch wide at 300 ppi"
$ws.Range("H21").Value = "44 inch wide at 300 ppi"
$ws.Range("H22").Value = "44 inch wide at 300 ppi (26 inches high"

# --- Numeric values
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = 20000
$ws.Range("C20").Value = 4500

$ws.Range("A21").Value = 2
$ws.Range("C21").Value = 13200

$ws.Range("B22").Value = 13200
$ws.Range("C22").Value = 13200

# --- Formulas
$ws.Range("F20").Formula = "=+B20*C20"
$ws.Range("F21").Formula = "=+`$F20"
$ws.Range("B21").Formula = "=+F21/C21"
$ws.Range("L22").Formula = "=13200/300"
$ws.Range("A26").Formula = "=+40*300"

# --- Column width for column C (~11.71 characters wide)
$ws.Columns.Item(3).ColumnWidth = 11.7109375

# --- Selection / active cell ends up on A27 after data entry
$ws.Range("A27").Select() | Out-Null

# --- Window position (best effort; matches author's saved view)
$excel.ActiveWindow.Left = 28800
$excel.ActiveWindow.Top = 525
